$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for 3 more EC (estado de cuenta) rows. Inserting whole rows
#    at 31:33 pushes the trailing "firma representante legal" footer block
#    (old rows 35:36) down to 38:39 while keeping its formatting/merged
#    cells intact.
# ---------------------------------------------------------------------------
$ws.Rows("31:33").Insert()

# Row 33 becomes the new "last" data row, so it needs the special bottom
# border style that used to live on row 30. Grab it before row 30 is
# overwritten below.
$ws.Range("B30:J30").Copy($ws.Range("B33:J33"))

# Row 30 is no longer last, so it becomes a normal data row (style copied
# from a normal row, e.g. row 29). Rows 31/32 need the same normal style.
$ws.Range("B29:J29").Copy($ws.Range("B30:J30"))
$ws.Range("B29:J29").Copy($ws.Range("B31:J31"))
$ws.Range("B29:J29").Copy($ws.Range("B32:J32"))

# ---------------------------------------------------------------------------
# 2) Replace the EC detail table (rows 16-33) with the new data set.
# ---------------------------------------------------------------------------
$data = @(
    @("CC", "1128044380",      "OMAR ALONSO VARGAS MORELO",        "2206", 58579, 1464490),
    @("CC", "1051888166",      "JORGE LUIS LAMADRID MEDINA",       "2009", 17556, 689455),
    @("CC", "1047482816",      "DAMIAN PEREZ ACEVEDO",             "2108", 36341, 1200000),
    @("CC", "73008944",        "ELKIN RAFAEL SIERRA CARO",         "1802", 28290, 848714),
    @("CC", "79417905",        "JUAN ALBERTO HOYOS CUARTAS",       "1905", 33125, 908526),
    @("CC", "1048460036",      "ROSANGELA CAROLINA ROMERO BURGOS", "2012", 9363,  877803),
    @("CC", "93355255",        "LEONEL TORRES",                    "1912", 33125, 828116),
    @("CC", "20090288",        "CARLOS ANDRES BLANCO TUIRAN",      "1709", 29509, 781242),
    @("CC", "1047409424",      "JULIETT PAOLA ANGULO BEJARANO",    "2003", 35112, 828116),
    @("CC", "1047409424",      "JULIETT PAOLA ANGULO BEJARANO",    "2002", 35112, 828116),
    @("CC", "1102839947",      "KARINA MARIA BELLO GOMEZCACERES",  "2102", 13325, 908526),
    @("PE", "963298814051990", "STEPHANYE PATRICIA PERDOMO SAER",  "1912", 37276, 931889),
    @("CC", "20255153",        "ROSANGELA CAROLINA ROMERO BURGOS", "2004", 35112, 877803),
    @("CC", "20255153",        "ROSANGELA CAROLINA ROMERO BURGOS", "2003", 35112, 877803),
    @("CC", "20255153",        "ROSANGELA CAROLINA ROMERO BURGOS", "2002", 35112, 877803),
    @("CC", "30578736",        "OMAIRA ESTHER GUERRERO FLOREZ",    "2403", 34666, 1300000),
    @("CC", "30578736",        "OMAIRA ESTHER GUERRERO FLOREZ",    "2402", 52000, 1300000),
    @("CC", "30578736",        "OMAIRA ESTHER GUERRERO FLOREZ",    "2401", 8667,  1300000)
)

$row = 16
foreach ($rec in $data) {
    $ws.Range("B$row").Value = $rec[0]
    $ws.Range("C$row").Value = $rec[1]
    $ws.Range("D$row").Value = $rec[2]
    $ws.Range("E$row").Value = $rec[3]
    $ws.Range("F$row").Value = $rec[4]
    $ws.Range("G$row").Value = $rec[5]
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3) Update the summary header figures.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 567382
$ws.Range("C13").Value = 13
$ws.Range("F13").Value = 15
